# [EI-979] [Rollback] Change "Then_Question"/"Else_Question" header labels
# back to "Then_Goto"/"Else_Goto" in the survey.xlsx data dictionary (Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I1").Value = "Then_Goto"
$ws.Range("J1").Value = "Else_Goto"
